# Update to final (!) version of RACGP evidence table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Refreshed evidence-table figures (higher-precision recalculated values).
$ws.Range("D3").Value = 6.2685227799999996

$ws.Range("C4").Value = 5.2743417990000001
$ws.Range("D4").Value = 9.4110842520000002

$ws.Range("C5").Value = 7.9470410280000001
$ws.Range("D5").Value = 12.822451210000001

$ws.Range("C6").Value = 9.8308525109999998
$ws.Range("D6").Value = 15.010274730000001

$ws.Range("C7").Value = 6.4672902390000004
$ws.Range("E7").Value = 15.294215489999999
$ws.Range("F7").Value = 19.048720880000001

# Cost for "Duration" row was recomputed to a genuinely new figure.
$ws.Range("C8").Value = 3.923
$ws.Range("D8").Value = 7.4831314799999999
$ws.Range("E8").Value = 13.182373630000001

$ws.Range("C9").Value = 4.0721779700000003
$ws.Range("D9").Value = 6.6089478929999999

$ws.Range("C10").Value = 5.8733768309999999
$ws.Range("D10").Value = 13.191412809999999
$ws.Range("E10").Value = 17.647624629999999

# Leave the cursor where the author left it when they saved.
$ws.Range("A12").Select()
